# Fix #NUM! errors in column C (Saudi Arabia) for rows 2-6 on the
# "split_few_bars_nb0" sheet by replacing them with the correct computed
# numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.208838632887556
$ws.Range("C3").Value = 0.178857494432493
$ws.Range("C4").Value = 0.187422136325044
$ws.Range("C5").Value = 0.188547853061829
$ws.Range("C6").Value = 0.236333883293078
